$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.314.39"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "2.605.57"
$ws.Range("E3").Value = "  +9.52%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.49"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.40"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.601"
$ws.Range("E7").Value = "  +5.65%  "
$ws.Range("E9").Value = "  +12.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.99"
$ws.Range("E10").Value = "  +12.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  +5.74%  "
$ws.Range("E12").Value = "  +13.82%  "
$ws.Range("D13").Value = "3.003.19"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "2.630.05"
$ws.Range("E15").Value = "  +10.69%  "
$ws.Range("E16").Value = "  +10.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.89"
$ws.Range("E17").Value = "  +8.93%  "
$ws.Range("D18").Value = "46.463.30"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +4.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  +10.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.12"
$ws.Range("E22").Value = "  +5.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.50"
$ws.Range("E23").Value = "  +4.55%  "
$ws.Range("E24").Value = "  +7.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  +16.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.21"
$ws.Range("E26").Value = "  +33.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +6.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.06"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("E31").Value = "  +11.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.72"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("E33").Value = "  +20.18%  "
$ws.Range("E34").Value = "  +4.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0835"
$ws.Range("E35").Value = "  +7.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "151.69"
$ws.Range("E36").Value = "  +3.19%  "
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("E38").Value = "  +5.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.18"
$ws.Range("E39").Value = "  +7.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.63"
$ws.Range("E40").Value = "  +5.00%  "
$ws.Range("E41").Value = "  +11.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0323"
$ws.Range("E42").Value = "  +7.73%  "
$ws.Range("E43").Value = "  +6.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.35"
$ws.Range("E44").Value = "  +36.33%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.03"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.79"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.53"
$ws.Range("E48").Value = "  +11.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.13"
$ws.Range("E49").Value = "  +7.94%  "
$ws.Range("D50").Value = "2.859.97"
$ws.Range("E50").Value = "  +9.56%  "
$ws.Range("E51").Value = "  +7.72%  "
